$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 3-11 (columns D, L, M, N, O, P, Q, R, S, T) get their data
# permuted/reassigned to new weekly values. Apply the new values directly.

$newData = @{
    3  = @{ D = 44511; L = "Primera"; M = 45;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Los Andes"; S = 2800; T = 10 }
    4  = @{ D = 44511; L = "Primera"; M = 45;  N = 3200;  O = 3200;  P = 3200;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 320;  T = 10 }
    5  = @{ D = 44166; L = "Segunda"; M = 20;  N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";    R = "La Ligua";                S = 667;  T = 18 }
    6  = @{ D = 44466; L = "Primera"; M = 80;  N = 11000; O = 11000; P = 11000; Q = "`$/bandeja 5 kilos";  R = "La Ligua";                S = 2200; T = 5  }
    7  = @{ D = 44515; L = "Primera"; M = 80;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Los Andes"; S = 2800; T = 10 }
    8  = @{ D = 44496; L = "Primera"; M = 55;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2800; T = 10 }
    9  = @{ D = 44488; L = "Primera"; M = 100; N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 5 kilos";  R = "La Ligua";                S = 2400; T = 5  }
    10 = @{ D = 44483; L = "Primera"; M = 35;  N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 5 kilos";  R = "Provincia de Quillota";  S = 2000; T = 5  }
    11 = @{ D = 44519; L = "Primera"; M = 30;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2800; T = 10 }
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
